# The legacy GSC export gained a new day's data at the bottom and the
# oldest day (2025-08-24, the first data row) is dropped from the report.
# Deleting that row shifts every remaining row up by one, which is exactly
# what the authoritative diff shows (row 91 disappears, every C-column
# value moves up one row, and the now-unused "2025-08-24" shared string is
# removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the first data row (row 2, "2025-08-24"); Excel shifts the rest
# of the data up automatically, updating the sheet dimension and the
# shared-string table for us.
$ws.Rows.Item(2).Delete()
